$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new booking row (row 4) with data for Vesna Todoric
$ws.Range("A4").Value = "Vienna"
$ws.Range("B4").Value = "Belgrade"
$ws.Range("C4").Value = "'2022-03-05"
$ws.Range("D4").Value = "'2022-03-08"
$ws.Range("E4").Value = "vesna92@test.com"
$ws.Range("F4").Value = "Vesna"
$ws.Range("G4").Value = "Todoric"
$ws.Range("H4").Value = "Female"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "63522526"

# Add a hyperlink for the new email address, same as the existing one on E3
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:vesna92@test.com") | Out-Null

# Update the active selection to reflect the appropriate selector
$ws.Range("I6").Select() | Out-Null
